# Summary.xlsx (HW_2) — add a new "Proximity Search" model row to each of the
# four result tables on Sheet1 (STEMMED / STEMMED+With-HEAD / NON-STEMMED /
# NON-STEMMED+With-HEAD), right after the existing "Unigram LM with Laplace
# smoothing" row in each table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row right after each table's last data row (Unigram LM).
# Processed top-to-bottom; since each Insert() shifts everything below it
# down by one, the target row numbers below already account for the
# earlier inserts.
$ws.Rows("6:6").Insert()
$ws.Rows("12:12").Insert()
$ws.Rows("19:19").Insert()
$ws.Rows("25:25").Insert()

# STEMMED table — new Proximity Search row
$ws.Range("A6").Value = "Proximity Search"
$ws.Range("B6").Value = 0.2886
$ws.Range("C6").Value = 0.396
$ws.Range("D6").Value = 0.3293

# STEMMED / With HEAD – TEXT Indexed EC2 table — new Proximity Search row
$ws.Range("A12").Value = "Proximity Search"
$ws.Range("B12").Value = 0.2932
$ws.Range("C12").Value = 0.424
$ws.Range("D12").Value = 0.3293

# NON STEMMED table — new Proximity Search row
$ws.Range("A19").Value = "Proximity Search"
$ws.Range("B19").Value = 0.1544
$ws.Range("C19").Value = 0.26
$ws.Range("D19").Value = 0.2027

# NON STEMMED / With HEAD – TEXT Indexed EC2 table — new Proximity Search row
$ws.Range("A25").Value = "Proximity Search"
$ws.Range("B25").Value = 0.1605
$ws.Range("C25").Value = 0.268
$ws.Range("D25").Value = 0.2053

# Match the author's final cursor position recorded in the workbook view.
$ws.Range("D21").Select()
